$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to text format before writing so that
# numeric-looking values (e.g. "522.15") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '60.443.53'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '2.627.17'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '522.15'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('D6').Value = '151.87'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '0.576'
$ws.Range('E8').Value = '  -3.52%  '
$ws.Range('D9').Value = '6.44'
$ws.Range('E9').Value = '  -3.77%  '
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '3.088.14'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('D14').Value = '60.483.26'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '21.59'
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = '2.624.10'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('D19').Value = '348.49'
$ws.Range('E19').Value = '  -2.75%  '
$ws.Range('D20').Value = '10.49'
$ws.Range('E20').Value = '  -1.62%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '60.90'
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').Value = '0.165'
$ws.Range('E25').Value = '  -0.79%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('D27').Value = '0.0₃0840'
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').Value = '7.19'
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '6.08'
$ws.Range('E30').Value = '  +2.45%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = '19.14'
$ws.Range('E32').Value = '  -1.76%  '
$ws.Range('D33').Value = '149.84'
$ws.Range('E33').Value = '  -0.68%  '
$ws.Range('D34').Value = '4.03'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -1.98%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '0.889'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').Value = '0.891'
$ws.Range('E37').Value = '  +4.75%  '
$ws.Range('D38').Value = '36.50'
$ws.Range('E38').Value = '  +0.45%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Value = '1.46'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = '298.87'
$ws.Range('E40').Value = '  +2.07%  '
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('D42').Value = '0.630'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '0.0554'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = '19.76'
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').Value = '0.0237'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('E48').Value = '  -3.35%  '
$ws.Range('D49').Value = '10.36'
$ws.Range('E49').Value = '  +0.56%  '
$ws.Range('D50').Value = '18.96'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').Value = '1.972.25'
$ws.Range('E51').Value = '  -1.02%  '

# Restore the original (default) style on the Price column so no stray
# number-format / style attribute is left behind on the cells.
$ws.Range("D2:D51").Style = "Normal"
